$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells stay text-formatted (matching the source data, which stores
# prices/percentages as literal strings, not numbers) so values like "18.10"
# keep their trailing zero instead of being auto-coerced to 18.1 by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.651.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.574.96"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.76"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.86"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.96"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.035.32"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.554.37"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000144"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.562.62"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.19"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "337.58"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.64"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.72%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.97"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.701.94"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.54%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.47"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.86"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.21"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0810"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "460.55"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +10.29%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.45%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.89"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.47"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.68"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "157.65"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.71"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.02"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.626"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0534"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.10"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.76%  "
